$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2..67).
# The automatic update bumps that date by one day (45181 -> 45182, i.e.
# 2023-09-12 -> 2023-09-13) for every row, leaving everything else untouched.
for ($row = 2; $row -le 67; $row++) {
    $ws.Cells.Item($row, 3).Value = 45182
}
